$wb = $excel.ActiveWorkbook

# --- Data edit -------------------------------------------------------
# The "symptoms" column (E) on the "clinical values" sheet used the value
# "none" for patients without symptoms. The commit replaces that value with
# "asymptomatic" everywhere it's used (rows 2,3,5,6,7,8,9,10,12 - the rows
# whose symptoms cell literally read "none"; rows 1 (header), 4 and 11 use
# other values and stay untouched).
$wsClinical = $wb.Worksheets.Item("clinical values")
$wsClinical.Range("E2").Value = "asymptomatic"
$wsClinical.Range("E3").Value = "asymptomatic"
$wsClinical.Range("E5").Value = "asymptomatic"
$wsClinical.Range("E6").Value = "asymptomatic"
$wsClinical.Range("E7").Value = "asymptomatic"
$wsClinical.Range("E8").Value = "asymptomatic"
$wsClinical.Range("E9").Value = "asymptomatic"
$wsClinical.Range("E10").Value = "asymptomatic"
$wsClinical.Range("E12").Value = "asymptomatic"

# --- View / active-tab changes ----------------------------------------
# The workbook now opens with the "clinical values" sheet (3rd tab, index 2)
# active/selected, with cell E12 selected - instead of the "Family" sheet
# with G1 selected.
$wsClinical.Activate()
$wsClinical.Range("E12").Select()
